$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 309, shifting existing rows 309:328 down to 310:329
$ws.Rows(309).Insert()

# Fill the new row 309 with the weekly record, matching the surrounding rows' pattern
$ws.Range("A309").Value = 4
$ws.Range("B309").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C309").Value = "Los Lagos"
$ws.Range("D309").Value = 44826
$ws.Range("D309").NumberFormat = $ws.Range("D310").NumberFormat
$ws.Range("E309").Value = 10
$ws.Range("F309").Value = 100112003
$ws.Range("G309").Value = "Ajo"
$ws.Range("H309").Value = "Chino"
$ws.Range("I309").Value = "Primera"
$ws.Range("J309").Value = 100
$ws.Range("K309").Value = 23000
$ws.Range("L309").Value = 23000
$ws.Range("M309").Value = 23000
$ws.Range("N309").Value = "$/caja 10 kilos"
$ws.Range("O309").Value = "China"
$ws.Range("P309").Value = 2300
$ws.Range("Q309").Value = 10
$ws.Range("R309").Value = "Hortaliza"
